$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 with the master server's data.
# Order matters: new shared strings are appended in first-write order, and
# the target file expects them as 000106001(7), 127.0.0.1(8), MasterServer_1(9).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "000106001"

$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "MasterServer_1"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "MasterServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 2001

# Move/collapse the active selection to H3.
$ws.Range("H3").Select()
